$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 7 after the existing row 6, shifting nothing below it
# (row 6 is currently the last data row), and copy row 6's original values
# into it before updating row 6 with the new data.
$ws.Rows.Item(7).Insert()

$srcRange = $ws.Range("A6:T6")
$dstRange = $ws.Range("A7:T7")
$srcRange.Copy()
$dstRange.PasteSpecial(-4104)  # xlPasteAll

# Now update row 6 with the new values from the diff
$ws.Range("D6").Value = 44438
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100

$ws.Range("A1").Select()
